$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new time log entry: 2025-06-14 (serial 45822), 14:30 - 15:30, "CS Introduction Lecture 18"
# Carry over the date formatting (s="1", mm/dd/yyyy) from the row above instead
# of re-creating a style, then overwrite with the new date's value.
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A24").Value = 45822

$ws.Range("B24").Value = 14
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 30
$ws.Range("F24").Value = "CS Introduction Lecture 18"

# Update view: scroll so row 23 is at top, select E24 as active cell
$ws.Range("E24").Select()
$excel.ActiveWindow.ScrollRow = 23
